$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update _id and arrival-time-like date (latest_dose_date) only
$ws.Range("A2").Value = "612c446e21b1190043d71b5d"
$ws.Range("M2").Value = "Tue Sep 21 2021 04:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 3
$ws.Range("A3").Value = "612c44dc82a6e7f8ae653cf4"
$ws.Range("C3").Value = "MOHIT 2"
$ws.Range("D3").Value = "f20200047@pilani.bits-pilani.ac.in"
$ws.Range("M3").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 4
$ws.Range("A4").Value = "612c44f082a6e7f8ae653cf6"
$ws.Range("C4").Value = "MOHIT 3"
$ws.Range("D4").Value = "f20190047@pilani.bits-pilani.ac.in"
$ws.Range("M4").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 5
$ws.Range("A5").Value = "612c450282a6e7f8ae653cf8"
$ws.Range("C5").Value = "MOHIT 4"
$ws.Range("D5").Value = "f20180047@pilani.bits-pilani.ac.in"
$ws.Range("M5").Value = "Thu Jan 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 6
$ws.Range("A6").Value = "612c451c82a6e7f8ae653cfa"
$ws.Range("C6").Value = "MOHIT 5"
$ws.Range("D6").Value = "f20190000@pilani.bits-pilani.ac.in"
$ws.Range("M6").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 7
$ws.Range("A7").Value = "612c453082a6e7f8ae653cfc"
$ws.Range("C7").Value = "MOHIT 6"
$ws.Range("D7").Value = "f20180027@pilani.bits-pilani.ac.in"
$ws.Range("M7").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 8
$ws.Range("A8").Value = "612c454782a6e7f8ae653cfe"
$ws.Range("C8").Value = "MOHIT 7"
$ws.Range("D8").Value = "f20200037@pilani.bits-pilani.ac.in"
$ws.Range("M8").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 9
$ws.Range("A9").Value = "612c455582a6e7f8ae653d00"
$ws.Range("C9").Value = "MOHIT 8"
$ws.Range("D9").Value = "f20200041@pilani.bits-pilani.ac.in"
$ws.Range("M9").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 10
$ws.Range("A10").Value = "612c456882a6e7f8ae653d02"
$ws.Range("C10").Value = "MOHIT 9"
$ws.Range("D10").Value = "f20190042@pilani.bits-pilani.ac.in"
$ws.Range("M10").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 11
$ws.Range("A11").Value = "612c457b82a6e7f8ae653d04"
$ws.Range("C11").Value = "MOHIT 20"
$ws.Range("D11").Value = "f20180147@pilani.bits-pilani.ac.in"
$ws.Range("M11").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"

# Row 12: was AKSHAT KUMAR row with a different pic/pdf/consent structure.
# Replace entirely with the new MOHIT 11 record that matches the standard layout.
$ws.Range("A12").Value = "612c458f82a6e7f8ae653d06"
$ws.Range("B12").Value = "https://lh3.googleusercontent.com/a-/AOh14Gj1ww5UEswYptQPWEoVEaPYRkThY6c5A9AJQmVd=s96-c"
$ws.Range("C12").Value = "MOHIT 11"
$ws.Range("D12").Value = "f20191047@pilani.bits-pilani.ac.in"
$ws.Range("F12").Value = "FAILED"
$ws.Range("I12").Value = "media/pdf/f20200048@pilani.bits-pilani.ac.in.pdf"
$ws.Range("J12").Value = "media/consent/f20200048@pilani.bits-pilani.ac.in.pdf"
$ws.Range("M12").Value = "Tue Sep 21 2021 16:30:00 GMT+0000 (Coordinated Universal Time)"
